# database/industries/siman/sehegmat/official/yearly.xlsx
# "update database and change read_price algorithm"
#
# The twelve-month reporting window rolls forward one year: the
# 1396/10 column is dropped, 1397/10..1400/10 shift one column to the
# left (E<-F<-G<-H), and a brand-new 1401/10 column lands in I.
# Every data row under the two tables gets the same left-shift, with
# fresh 1401/10 figures (and, where the source changed, a fresh
# 1400/10 figure too) filling the vacated H/I cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- period headers (row 8 and row 24) ---
$periods = @("دوازده ماهه منتهی به 1397/10", "دوازده ماهه منتهی به 1398/10", "دوازده ماهه منتهی به 1399/10", "دوازده ماهه منتهی به 1400/10", "دوازده ماهه منتهی به 1401/10")
$cols = @("E", "F", "G", "H", "I")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $periods[$i]
    $ws.Range($cols[$i] + "24").Value = $periods[$i]
}

# --- table 1 (هزینه های عمومی و اداری), rows 10-20 ---
$table1 = @{
    10 = @(48037, 380069, 335198, 540757, 1084979)
    13 = @(1714, 369, 0, 0, 0)
    15 = @(1201, 405, 362, 0, 0)
    16 = @(809, 507, 634, 3518, 406)
    17 = @(55896, 81282, 69980, 74076, 121359)
    18 = @(0, 18047, 6092, 25788, 1119)
    19 = @(40303, 62622, 158405, 124054, 163135)
    20 = @(147960, 543301, 570671, 768193, 1370998)
}

foreach ($row in $table1.Keys) {
    $values = $table1[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}

# --- table 2 (تعداد پرسنل), rows 26-27 ---
$table2 = @{
    26 = @(337, 68, 58, 67, 65)
    27 = @(61, 389, 386, 674, 375)
}

foreach ($row in $table2.Keys) {
    $values = $table2[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}
